# The post at row 216 ("「どうやらこの好奇心旺盛なペンギンは…」" - the curious
# king penguin post) was removed from the source data. Delete that entire
# row from the worksheet; Excel shifts every row below it up by one,
# which matches the renumbering seen across the rest of the diff
# (old row 217 -> new row 216, ..., old row 316 -> new row 315) and the
# updated used-range dimension A1:C315.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(216).Delete()
